# Update PO_Forecast column (B) for rows 10 through 25 from 47 to 48
# This implements the "added 4wk low sales check" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 10; $r -le 25; $r++) {
    $ws.Cells.Item($r, 2).Value = 48
}
